$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values: switch the touched range to
# Text format while assigning so Excel does not coerce strings like "1.00" or
# "649.20" into numbers, then reset the style back to Normal afterwards so no
# stray style index (`s` attribute) is left on the cell.

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "69.990.58"
$ws.Range("E2").Value = "  +0.61%  "
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "3.698.23"
$ws.Range("E3").Value = "  +0.08%  "
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "649.20"
$ws.Range("E5").Value = "  -3.98%  "
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "161.94"
$ws.Range("E6").Value = "  +0.49%  "
$rng.Style = "Normal"

$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$rng.Style = "Normal"

$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$ws.Range("E8").Value = "  +1.50%  "
$rng.Style = "Normal"

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$ws.Range("E9").Value = "  -0.60%  "
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  +1.34%  "
$rng.Style = "Normal"

$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$ws.Range("E11").Value = "  +0.99%  "
$rng.Style = "Normal"

$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$ws.Range("E12").Value = "  -0.18%  "
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "4.320.40"
$ws.Range("E13").Value = "  -0.02%  "
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "32.87"
$ws.Range("E14").Value = "  +0.97%  "
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "3.712.20"
$ws.Range("E15").Value = "  +0.21%  "
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "69.952.62"
$ws.Range("E16").Value = "  +0.60%  "
$rng.Style = "Normal"

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$ws.Range("E17").Value = "  +0.62%  "
$rng.Style = "Normal"

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = "  +0.47%  "
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = "  +1.09%  "
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = "10.50"
$ws.Range("E20").Value = "  +6.64%  "
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "472.35"
$ws.Range("E21").Value = "  +0.28%  "
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("E22").Value = "  +0.55%  "
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "80.19"
$ws.Range("E23").Value = "  -0.43%  "
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "3.842.25"
$ws.Range("E24").Value = "  -0.02%  "
$rng.Style = "Normal"

$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$ws.Range("E25").Value = "  +2.11%  "
$rng.Style = "Normal"

$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$rng.Style = "Normal"

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  +1.00%  "
$rng.Style = "Normal"

$rng = $ws.Range("D28:E28")
$rng.NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").Value = "  +0.98%  "
$rng.Style = "Normal"

$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$ws.Range("E29").Value = "  -1.54%  "
$rng.Style = "Normal"

$rng = $ws.Range("D31:E31")
$rng.NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  +0.51%  "
$rng.Style = "Normal"

$rng = $ws.Range("D32:E32")
$rng.NumberFormat = "@"
$ws.Range("D32").Value = "6.57"
$ws.Range("E32").Value = "  -0.32%  "
$rng.Style = "Normal"

$rng = $ws.Range("B33:E33")
$rng.NumberFormat = "@"
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.167"
$ws.Range("E33").Value = "  +3.33%  "
$rng.Style = "Normal"

$rng = $ws.Range("B34:E34")
$rng.NumberFormat = "@"
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.09%  "
$rng.Style = "Normal"

$rng = $ws.Range("D35:E35")
$rng.NumberFormat = "@"
$ws.Range("D35").Value = "26.85"
$ws.Range("E35").Value = "  -0.48%  "
$rng.Style = "Normal"

$rng = $ws.Range("D36:E36")
$rng.NumberFormat = "@"
$ws.Range("D36").Value = "3.693.03"
$ws.Range("E36").Value = "  +0.17%  "
$rng.Style = "Normal"

$rng = $ws.Range("D37:E37")
$rng.NumberFormat = "@"
$ws.Range("D37").Value = "8.46"
$ws.Range("E37").Value = "  +0.19%  "
$rng.Style = "Normal"

$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$ws.Range("E39").Value = "  -4.78%  "
$rng.Style = "Normal"

$rng = $ws.Range("B40:E40")
$rng.NumberFormat = "@"
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "179.89"
$ws.Range("E40").Value = "  +7.68%  "
$rng.Style = "Normal"

$rng = $ws.Range("B41:E41")
$rng.NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.25"
$ws.Range("E41").Value = "  +0.82%  "
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.06%  "
$rng.Style = "Normal"

$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$ws.Range("E43").Value = "  +0.43%  "
$rng.Style = "Normal"

$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$ws.Range("E44").Value = "  -1.04%  "
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("E45").Value = "  +4.39%  "
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "29.36"
$ws.Range("E46").Value = "  +4.18%  "
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "46.98"
$ws.Range("E47").Value = "  +0.80%  "
$rng.Style = "Normal"

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$ws.Range("E48").Value = "  -0.96%  "
$rng.Style = "Normal"

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = "  -1.84%  "
$rng.Style = "Normal"

$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$ws.Range("E50").Value = "  -0.02%  "
$rng.Style = "Normal"

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$ws.Range("E51").Value = "  -3.16%  "
$rng.Style = "Normal"
